$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.442.49"
$ws.Range("E2").Value = "  -3.11%  "

$ws.Range("D3").Value = "3.487.92"
$ws.Range("E3").Value = "  -0.70%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.80"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.77"
$ws.Range("E6").Value = "  -5.91%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.639"
$ws.Range("E7").Value = "  +4.31%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.632"
$ws.Range("E9").Value = "  -1.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("E10").Value = "  +2.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.63"
$ws.Range("E11").Value = "  -6.07%  "

$ws.Range("E12").Value = "  -1.93%  "

$ws.Range("E13").Value = "  -3.05%  "

$ws.Range("D14").Value = "4.051.21"
$ws.Range("E14").Value = "  -0.25%  "

$ws.Range("D15").Value = "3.494.28"
$ws.Range("E15").Value = "  -0.23%  "

$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.35"
$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("E18").Value = "  +1.20%  "

$ws.Range("D19").Value = "65.531.34"
$ws.Range("E19").Value = "  -4.11%  "

$ws.Range("E20").Value = "  -1.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "413.34"
$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("E22").Value = "  +2.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "85.94"
$ws.Range("E23").Value = "  +1.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.10"
$ws.Range("E24").Value = "  -2.81%  "

$ws.Range("E25").Value = "  +6.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.79"
$ws.Range("E26").Value = "  -8.41%  "

$ws.Range("E27").Value = "  -2.20%  "

$ws.Range("E28").Value = "  -3.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.03"
$ws.Range("E29").Value = "  +4.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.23"
$ws.Range("E30").Value = "  -1.17%  "

$ws.Range("E31").Value = "  -5.55%  "

$ws.Range("E32").Value = "  -11.59%  "

$ws.Range("E33").Value = "  -0.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("E34").Value = "  -1.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.45"
$ws.Range("E35").Value = "  -1.44%  "

$ws.Range("E36").Value = "  +10.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0789"
$ws.Range("E38").Value = "  -5.06%  "

$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.06"
$ws.Range("E39").Value = "  -5.44%  "

$ws.Range("D40").Value = "3.351.98"
$ws.Range("E40").Value = "  +9.63%  "

$ws.Range("E41").Value = "  -5.78%  "

$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("E43").Value = "  -3.90%  "

$ws.Range("E44").Value = "  -6.75%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.52"
$ws.Range("E45").Value = "  -9.56%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0414"
$ws.Range("E46").Value = "  -1.89%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.25"
$ws.Range("E47").Value = "  +0.42%  "

$ws.Range("E48").Value = "  -1.61%  "

$ws.Range("E49").Value = "  +1.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.42"
$ws.Range("E50").Value = "  -5.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.07"
$ws.Range("E51").Value = "  -1.80%  "
